{"js": "// The document contains a paragraph with bold text that reads\n// \"DOCX, DOC, PDF, HTML, XPS, R\" + a \"_GoBack\" bookmark + \"TF and TXT\"\n// split across two runs (an artifact of a prior edit/cursor position).\n// The fix merges them back into a single run \"DOCX, DOC, PDF, HTML, XPS,\n// RTF and TXT\" and removes the now-stale \"_GoBack\" bookmark.\n\nconst body = context.document.body;\n\n// Word's Find/search treats bookmarks as invisible, so the text is seen\n// as contiguous even though it's split across two runs around the\n// bookmark.\nconst results = body.search(\"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\", { matchCase: true });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  // Replacing the range's text with itself collapses the underlying runs\n  // (and the bookmark boundary sitting between them) into a single run,\n  // while keeping the existing (bold) formatting of the range.\n  target.insertText(\"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\", Word.InsertLocation.replace);\n}\n\n// Remove the stray \"_GoBack\" bookmark left over in the paragraph.\ncontext.document.deleteBookmark(\"_GoBack\");\n\nawait context.sync();\n", "ps1": "# The document contains a paragraph with bold text that reads\n# \"DOCX, DOC, PDF, HTML, XPS, R\" + a \"_GoBack\" bookmark + \"TF and TXT\"\n# split across two runs (an artifact of a prior edit/cursor position).\n# The fix merges them back into a single run \"DOCX, DOC, PDF, HTML, XPS,\n# RTF and TXT\" and removes the now-stale \"_GoBack\" bookmark.\n\n$d = $word.ActiveDocument\n\n# Word's Find treats bookmarks as invisible, so the text is found as\n# contiguous even though it's split across two runs around the bookmark.\n$range = $d.Content\n$find = $range.Find\n$find.Text = \"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\"\n$find.Replacement.Text = \"DOCX, DOC, PDF, HTML, XPS, RTF and TXT\"\n\n# Replace-in-place: collapses the split runs (and the bookmark boundary\n# sitting between them) back into a single run while keeping the\n# existing (bold) formatting of the matched text.\n$found = $find.Execute(\n    $find.Text, $false, $false, $false, $false, $false, $true, 1, $false,\n    $find.Replacement.Text, 2\n)\n\n# Remove the stray \"_GoBack\" bookmark left over in the paragraph.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks(\"_GoBack\").Delete()\n}\n"}
